$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the formatting of columns A and B (row 32's look) down to the four
# brand-new rows (33-36) before any values are written, so the new cells
# pick up the same cell style (shading/border) as the rest of the fuel list
# without creating duplicate style entries in styles.xml.
$ws.Range("A32").Copy()
$ws.Range("A33:A36").PasteSpecial(-4122)
$ws.Range("B32").Copy()
$ws.Range("B33:B36").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rewrite the fuel list (columns A and B, rows 3-36) with its final content.
# Four new fuels - "Biomass for industry", "Biomass for other sectors",
# "Coal for industry" and "Coal for other sectors" - are inserted right
# after their parent fuel (Biomass / Coal), so every row below shifts down
# by the corresponding amount; we simply write out the resulting values.
$ws.Range("A3").Value = "Biomass for industry"
$ws.Range("B3").Value = "BIO_IND"
$ws.Range("A4").Value = "Biomass for other sectors"
$ws.Range("B4").Value = "BIO_OTH"
$ws.Range("A5").Value = "Coal"
$ws.Range("B5").Value = "COA"
$ws.Range("A6").Value = "Coal for industry"
$ws.Range("B6").Value = "COA_IND"
$ws.Range("A7").Value = "Coal for other sectors"
$ws.Range("B7").Value = "COA_OTH"
$ws.Range("A8").Value = "Commercial electricity"
$ws.Range("B8").Value = "COM_ELC"
$ws.Range("A9").Value = "Commercial other energy"
$ws.Range("B9").Value = "COM_OTH"
$ws.Range("A10").Value = "diesel"
$ws.Range("B10").Value = "DSL"
$ws.Range("A11").Value = "Electricity for transmission"
$ws.Range("B11").Value = "ELC_001"
$ws.Range("A12").Value = "Electricity for final uses"
$ws.Range("B12").Value = "ELC_002"
$ws.Range("A13").Value = "Electricity for final uses - new connections"
$ws.Range("B13").Value = "ELC_003"
$ws.Range("A14").Value = "Geothermal"
$ws.Range("B14").Value = "GEO"
$ws.Range("A15").Value = "Gasoline"
$ws.Range("B15").Value = "GSL"
$ws.Range("A16").Value = "Hydro"
$ws.Range("B16").Value = "HYD"
$ws.Range("A17").Value = "Industrial energy uses"
$ws.Range("B17").Value = "IND_EN"
$ws.Range("A18").Value = "Kerosene"
$ws.Range("B18").Value = "KER"
$ws.Range("A19").Value = "LPG"
$ws.Range("B19").Value = "LPG"
$ws.Range("A20").Value = "Natural Gas"
$ws.Range("B20").Value = "NGS"
$ws.Range("A21").Value = "Oil"
$ws.Range("B21").Value = "OIL"
$ws.Range("A22").Value = "Other oil products"
$ws.Range("B22").Value = "OILPROD"
$ws.Range("A23").Value = "Other energy uses"
$ws.Range("B23").Value = "OTH_EN"
$ws.Range("A24").Value = "Residential cooling"
$ws.Range("B24").Value = "RES_COOL"
$ws.Range("A25").Value = "Residential cooling - new connections"
$ws.Range("B25").Value = "RES_COOLb"
$ws.Range("A26").Value = "Residential energy for cooking and water heating"
$ws.Range("B26").Value = "RES_CWH"
$ws.Range("A27").Value = "Residential electricity for appliances"
$ws.Range("B27").Value = "RES_ELC_APP"
$ws.Range("A28").Value = "Residential electricity for appliances - new connections"
$ws.Range("B28").Value = "RES_ELC_APPb"
$ws.Range("A29").Value = "Solar energy"
$ws.Range("B29").Value = "SOL"
$ws.Range("A30").Value = "Transportation freight aviation & navigation - in Mton-km"
$ws.Range("B30").Value = "TRA_AN_FREIGHT"
$ws.Range("A31").Value = "Transportation passenger aviation & navigation - in Mpassenger-km"
$ws.Range("B31").Value = "TRA_AN_PSNG"
$ws.Range("A32").Value = "Transportation railway freight - in Mton-km"
$ws.Range("B32").Value = "TRA_RLW_FREIGHT"
$ws.Range("A33").Value = "Transportation passenger railway - in Mpassenger - km"
$ws.Range("B33").Value = "TRA_RLW_PSNG"
$ws.Range("A34").Value = "Transportation road freight - in Mton-km"
$ws.Range("B34").Value = "TRA_ROAD_FREIGHT"
$ws.Range("A35").Value = "Transportation road passenger - in Mpassenger-km"
$ws.Range("B35").Value = "TRA_ROAD_PSNG"
$ws.Range("A36").Value = "Wind energy"
$ws.Range("B36").Value = "WND"

# Restore the selected cell that Excel recorded after the edits were made.
$null = $ws.Range("B8").Select()
